$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7830654382705688
$ws.Range("B1").Value = 1.470719695091248
$ws.Range("C1").Value = 5.644104957580566
$ws.Range("D1").Value = 3.154501914978027
$ws.Range("E1").Value = 1.487643957138062
